$d = $word.ActiveDocument

# Locate the "Bases de données" skills paragraph as our anchor.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "^Bases de don") {
        $anchorIndex = $i
        break
    }
}

# Insert a new "Langages" paragraph (without "r") right before "Bases de données".
$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphBefore()
$d.Paragraphs.Item($anchorIndex).Range.Text = "Langages : python, matlab, c, c++"
$anchorIndex = $anchorIndex + 1

# Insert a new "Data Science" paragraph right after "Bases de données".
$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphAfter()
$d.Paragraphs.Item($anchorIndex + 1).Range.Text = "Data Science : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn"

# Find the old "MLOps" / "ML/AI" / "Langages : r, ..." paragraphs.
$mlopsIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "^MLOps") {
        $mlopsIndex = $i
        break
    }
}

# Rename "MLOps" to "Machine Learning" (text/content unchanged otherwise).
$d.Paragraphs.Item($mlopsIndex).Range.Text = "Machine Learning : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit"

# Remove the old "ML/AI : ..." paragraph entirely.
$d.Paragraphs.Item($mlopsIndex + 1).Range.Delete()

# Remove the old "Langages : r, python, matlab, c, c++" paragraph entirely.
$d.Paragraphs.Item($mlopsIndex + 1).Range.Delete()
